$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 with new numeric / boolean values ---

# row 8 (extr1)
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11

# row 9 (extr2)
$ws.Cells.Item(9, 3).Value = 16

# row 10 (extr3)
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# row 11 (extr4)
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# row 12 (extr5)
$ws.Cells.Item(12, 3).Value = 10

# row 13 (extr6)
$ws.Cells.Item(13, 4).Value = 8

# row 14 (extr7)
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11

# row 15 (extr8) - values shift; new values for row 15
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# --- Append two new rows for line7 / line8 ---

# Copy formatting of column A (bold/border/center style) down to the new rows
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

# row 16 (line7)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "line7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# row 17 (line8)
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "line8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true
